$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.994.77"
$ws.Range("E2").Value = "  -4.91%  "
$ws.Range("D3").Value = "2.225.92"
$ws.Range("E3").Value = "  -5.80%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "318.39"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "99.46"
$ws.Range("E6").Value = "  -7.92%  "
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  -8.83%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "0.559"
$ws.Range("E9").Value = "  -8.62%  "
$ws.Range("D10").Value = "36.80"
$ws.Range("E10").Value = "  -9.78%  "
$ws.Range("D11").Value = "54.07"
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("D12").Value = "0.0829"
$ws.Range("E12").Value = "  -9.72%  "
$ws.Range("D13").Value = "7.65"
$ws.Range("E13").Value = "  -9.74%  "
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "2.566.10"
$ws.Range("E15").Value = "  -5.84%  "
$ws.Range("D16").Value = "0.862"
$ws.Range("E16").Value = "  -11.93%  "
$ws.Range("D17").Value = "14.36"
$ws.Range("E17").Value = "  -6.43%  "
$ws.Range("D18").Value = "2.220.22"
$ws.Range("E18").Value = "  -5.76%  "
$ws.Range("D19").Value = "42.927.95"
$ws.Range("E19").Value = "  -4.95%  "
$ws.Range("D20").Value = "14.42"
$ws.Range("E20").Value = "  -5.09%  "
$ws.Range("D21").Value = "0.0₃0963"
$ws.Range("E21").Value = "  -9.27%  "
$ws.Range("D22").Value = "6.52"
$ws.Range("E22").Value = "  -10.10%  "
$ws.Range("D23").Value = "65.15"
$ws.Range("E23").Value = "  -10.95%  "
$ws.Range("D24").Value = "3.15"
$ws.Range("E24").Value = "  -12.42%  "
$ws.Range("D25").Value = "236.32"
$ws.Range("E25").Value = "  -10.34%  "
$ws.Range("D26").Value = "2.15"
$ws.Range("E26").Value = "  -8.03%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "4.02"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").Value = "9.97"
$ws.Range("E30").Value = "  -10.53%  "
$ws.Range("D31").Value = "6.32"
$ws.Range("E31").Value = "  -14.49%  "
$ws.Range("D32").Value = "35.75"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "20.27"
$ws.Range("E33").Value = "  -9.42%  "
$ws.Range("D34").Value = "0.0870"
$ws.Range("E34").Value = "  -7.99%  "
$ws.Range("D35").Value = "153.51"
$ws.Range("E35").Value = "  -9.19%  "
$ws.Range("D36").Value = "2.66"
$ws.Range("E36").Value = "  -6.30%  "
$ws.Range("D37").Value = "3.15"
$ws.Range("E37").Value = "  +6.01%  "
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  -7.64%  "
$ws.Range("D40").Value = "4.42"
$ws.Range("E40").Value = "  -5.60%  "
$ws.Range("D41").Value = "0.103"
$ws.Range("E41").Value = "  -11.14%  "
$ws.Range("D42").Value = "3.64"
$ws.Range("E42").Value = "  -7.92%  "
$ws.Range("E43").Value = "  -8.70%  "
$ws.Range("D44").Value = "13.46"
$ws.Range("E44").Value = "  +4.73%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "1.742.06"
$ws.Range("E46").Value = "  -6.88%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "85.12"
$ws.Range("E47").Value = "  -15.56%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.204"
$ws.Range("E48").Value = "  -10.07%  "
$ws.Range("D49").Value = "5.29"
$ws.Range("E49").Value = "  -10.73%  "
$ws.Range("D50").Value = "75.02"
$ws.Range("D51").Value = "8.67"
$ws.Range("E51").Value = "  -5.41%  "
